$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.243.62"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -1.03%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.662.88"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -1.06%  "
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.44%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "219.15"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +1.13%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5228"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -1.80%  "
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.39%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2672"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.37%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06352"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.74%  "
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.66%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07726"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -1.00%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.665.98"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -0.91%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.436"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -1.42%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "1.890.67"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.04%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5479"
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "0.0₅8225"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.30%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "65.07"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.92%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "26.259.73"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -1.11%  "
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  +0.55%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.660"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -2.07%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "195.78"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.60%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -1.85%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.095"
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -4.00%  "
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  +0.57%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "139.23"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -2.88%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1245"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -2.97%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.237"
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -2.67%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -0.84%  "
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -0.82%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05984"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -3.29%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.284"
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.82%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.622"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.313"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -4.12%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.634"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -3.25%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9819"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -2.61%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.426"
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  +0.08%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.783"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -0.26%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5907"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  +3.01%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01599"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -2.34%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "6.000"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -0.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8588"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.18%  "
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  +0.45%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.031.87"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -4.06%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.00"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.804.78"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -1.22%  "
$ws.Range("B46").Value = "BabyDogeCoin"
$ws.Range("C46").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.0₈113"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  +3.36%  "
$ws.Range("B47").Value = "Aave"
$ws.Range("C47").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "57.53"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  +0.89%  "
$ws.Range("B48").Value = "Frax"
$ws.Range("C48").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.011"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +1.06%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "8.077"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.75%  "
$ws.Range("B50").Value = "Cronos"
$ws.Range("C50").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05187"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -0.41%  "
$ws.Range("B51").Value = "RenderToken"
$ws.Range("C51").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.473"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  +0.26%  "
